$d = $word.ActiveDocument

$r1 = $d.Content
$found1 = $r1.Find.Execute("                            Sep 2", $false, $false, $false, $false, $false, $true, 1, $false, "                           Aug 2", 2)
if (-not $found1) { Write-Output "WARNING: site 1 search text not found" }

$r2 = $d.Content
$found2 = $r2.Find.Execute("`t`t`t`t`t         `t`t`t Jan ", $false, $false, $false, $false, $false, $true, 1, $false, "`t`t`t`t`t         `t`t           Nov ", 2)
if (-not $found2) { Write-Output "WARNING: site 2 search text not found" }

$r3 = $d.Content
$found3 = $r3.Find.Execute("              Jun 2012", $false, $false, $false, $false, $false, $true, 1, $false, "               Jun 2012", 2)
if (-not $found3) { Write-Output "WARNING: site 3 search text not found" }

$r4 = $d.Content
$found4 = $r4.Find.Execute("University of Windsor`t`t`t`t`t`t`t`t`t`t`t     2011", $false, $false, $false, $false, $false, $true, 1, $false, "University of Windsor`t`t`t`t`t`t`t`t`t`t`t2011", 2)
if (-not $found4) { Write-Output "WARNING: site 4 search text not found" }

$r5 = $d.Content
$found5 = $r5.Find.Execute("University of Windsor`t`t`t`t`t`t`t`t`t     2009", $false, $false, $false, $false, $false, $true, 1, $false, "University of Windsor`t`t`t`t`t`t`t`t`t2009", 2)
if (-not $found5) { Write-Output "WARNING: site 5 search text not found" }
